$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 544.8333
$ws.Range("I39").Value = 544.8333
$ws.Range("K39").Value = 1634.4999
$ws.Range("M39").Value = -1338.4999
# Row 64
$ws.Range("H64").Value = 3308.5
$ws.Range("J64").Value = 3350
$ws.Range("L64").Value = 3350
$ws.Range("N64").Value = -3846
# Row 67
$ws.Range("H67").Value = 3308.5
$ws.Range("J67").Value = 3350
$ws.Range("L67").Value = 3350
$ws.Range("N67").Value = -5066
# Row 125
$ws.Range("H125").Value = 1781.125
$ws.Range("J125").Value = 1975
$ws.Range("L125").Value = 17775
$ws.Range("N125").Value = -22695
# Row 132
$ws.Range("H132").Value = 886.675
$ws.Range("I132").Value = 886.675
$ws.Range("K132").Value = 2660.025
$ws.Range("M132").Value = -130.0249999999996
# Row 137
$ws.Range("H137").Value = 1869.6
$ws.Range("I137").Value = 1472.8
$ws.Range("K137").Value = 4418.4
$ws.Range("M137").Value = -1868.4
# Row 138
$ws.Range("H138").Value = 4177.706
$ws.Range("J138").Value = 4507.136
$ws.Range("L138").Value = 13521.408
$ws.Range("N138").Value = -23801.408
# Row 141
$ws.Range("H141").Value = 8999.75
$ws.Range("I141").Value = 7999
$ws.Range("J141").Value = 9333.333000000001
$ws.Range("K141").Value = 23997
$ws.Range("L141").Value = 27999.999
$ws.Range("M141").Value = -18817
$ws.Range("N141").Value = -38359.999

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2129.6667
$ws.Range("I2").Value = 2195
$ws.Range("J2").Value = 1999
$ws.Range("K2").Value = 2195
$ws.Range("L2").Value = 1999
$ws.Range("M2").Value = -2082
$ws.Range("N2").Value = -2225
# Row 32
$ws.Range("H32").Value = 5891.875
$ws.Range("I32").Value = 5671.1665
$ws.Range("J32").Value = 14499.5
$ws.Range("K32").Value = 5671.1665
$ws.Range("L32").Value = 14499.5
$ws.Range("M32").Value = -5384.1665
$ws.Range("N32").Value = -15073.5
# Row 45
$ws.Range("H45").Value = 2602.7144
$ws.Range("I45").Value = 2602.7144
$ws.Range("K45").Value = 2602.7144
$ws.Range("M45").Value = -2225.7144
# Row 61
$ws.Range("H61").Value = 1805.4166
$ws.Range("I61").Value = 1577
$ws.Range("J61").Value = 2947.5
$ws.Range("K61").Value = 1577
$ws.Range("L61").Value = 2947.5
$ws.Range("M61").Value = -1365
$ws.Range("N61").Value = -3371.5
# Row 74
$ws.Range("H74").Value = 2164.2632
$ws.Range("I74").Value = 1370.8462
$ws.Range("K74").Value = 1370.8462
$ws.Range("M74").Value = -496.8462
# Row 77
$ws.Range("H77").Value = 2164.2632
$ws.Range("I77").Value = 1370.8462
$ws.Range("K77").Value = 6854.231
$ws.Range("M77").Value = -2486.231
# Row 116
$ws.Range("H116").Value = 2129.6667
$ws.Range("I116").Value = 2195
$ws.Range("J116").Value = 1999
$ws.Range("K116").Value = 2195
$ws.Range("L116").Value = 1999
$ws.Range("M116").Value = 99
$ws.Range("N116").Value = -6587
# Row 122
$ws.Range("H122").Value = 8166.3335
$ws.Range("I122").Value = 8166.3335
$ws.Range("K122").Value = 24499.0005
$ws.Range("M122").Value = -22049.0005
# Row 132
$ws.Range("H132").Value = 3033.7273
$ws.Range("I132").Value = 2124.7856
$ws.Range("J132").Value = 4624.375
$ws.Range("K132").Value = 6374.3568
$ws.Range("L132").Value = 13873.125
$ws.Range("M132").Value = -3844.3568
$ws.Range("N132").Value = -18933.125
# Row 135
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
# Row 136
$ws.Range("H136").Value = 1805.4166
$ws.Range("I136").Value = 1577
$ws.Range("J136").Value = 2947.5
$ws.Range("K136").Value = 4731
$ws.Range("L136").Value = 8842.5
$ws.Range("M136").Value = -2181
$ws.Range("N136").Value = -13942.5
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2129.6667
$ws.Range("I3").Value = 2195
$ws.Range("J3").Value = 1999
$ws.Range("K3").Value = 2195
$ws.Range("L3").Value = 1999
$ws.Range("M3").Value = -2081
$ws.Range("N3").Value = -2227
# Row 134
$ws.Range("H134").Value = 3336.0557
$ws.Range("I134").Value = 2984.7
$ws.Range("J134").Value = 3775.25
$ws.Range("K134").Value = 8954.099999999999
$ws.Range("L134").Value = 11325.75
$ws.Range("M134").Value = -6419.099999999999
$ws.Range("N134").Value = -16395.75
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1921
$ws.Range("I16").Value = 1842
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1842
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1555
$ws.Range("N16").Value = -2574
# Row 58
$ws.Range("H58").Value = 2952.5
$ws.Range("J58").Value = 2979.75
$ws.Range("L58").Value = 2979.75
$ws.Range("N58").Value = -3385.75
# Row 113
$ws.Range("H113").Value = 1921
$ws.Range("I113").Value = 1842
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1842
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 328
$ws.Range("N113").Value = -6340
# Row 132
$ws.Range("H132").Value = 3651.077
$ws.Range("I132").Value = 3384.3684
$ws.Range("J132").Value = 4375
$ws.Range("K132").Value = 10153.1052
$ws.Range("L132").Value = 13125
$ws.Range("M132").Value = -7623.1052
$ws.Range("N132").Value = -18185
# Row 136
$ws.Range("H136").Value = 2952.5
$ws.Range("J136").Value = 2979.75
$ws.Range("L136").Value = 8939.25
$ws.Range("N136").Value = -14039.25
# Row 141
$ws.Range("H141").Value = 50369.43
$ws.Range("J141").Value = 50369.43
$ws.Range("L141").Value = 50369.43
$ws.Range("N141").Value = -60729.43

$ws = $wb.Worksheets.Item("CUL")
# Row 94
$ws.Range("H94").Value = 2777.8
$ws.Range("I94").Value = 1989
$ws.Range("K94").Value = 5967
$ws.Range("M94").Value = -5291
# Row 109
$ws.Range("H109").Value = 4222.1875
$ws.Range("J109").Value = 4642.857
$ws.Range("L109").Value = 13928.571
$ws.Range("N109").Value = -16008.571
# Row 131
$ws.Range("H131").Value = 1668.8572
$ws.Range("J131").Value = 2165.3333
$ws.Range("L131").Value = 6495.999899999999
$ws.Range("N131").Value = -16575.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 495.12
$ws.Range("I97").Value = 485
$ws.Range("K97").Value = 485
$ws.Range("M97").Value = 11
# Row 102
$ws.Range("H102").Value = 1249.5
$ws.Range("J102").Value = 999
$ws.Range("L102").Value = 999
$ws.Range("N102").Value = -4243
# Row 122
$ws.Range("H122").Value = 2054.625
$ws.Range("I122").Value = 2062.4285
$ws.Range("K122").Value = 6187.2855
$ws.Range("M122").Value = -3737.2855
# Row 126
$ws.Range("H126").Value = 1480.4
$ws.Range("I126").Value = 1478
$ws.Range("J126").Value = 1484
$ws.Range("K126").Value = 4434
$ws.Range("L126").Value = 4452
$ws.Range("M126").Value = -1964
$ws.Range("N126").Value = -9392
# Row 132
$ws.Range("H132").Value = 4674.6
$ws.Range("I132").Value = 4647.4546
$ws.Range("J132").Value = 4749.25
$ws.Range("K132").Value = 13942.3638
$ws.Range("L132").Value = 14247.75
$ws.Range("M132").Value = -11412.3638
$ws.Range("N132").Value = -19307.75
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2151.8
$ws.Range("I126").Value = 2091.2144
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6273.6432
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -3803.6432
$ws.Range("N126").Value = -13940
# Row 132
$ws.Range("H132").Value = 3787.6
$ws.Range("I132").Value = 2981.6667
$ws.Range("J132").Value = 4996.5
$ws.Range("K132").Value = 8945.000100000001
$ws.Range("L132").Value = 14989.5
$ws.Range("M132").Value = -6415.000100000001
$ws.Range("N132").Value = -20049.5
